# This workbook is being reverted from its "Solved" state (with formulas)
# back to a "Starter" state: the formula cells keep their number
# formatting/styles but lose their formulas and cached values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the Average/Final-Score formulas for every student row (F2:F10, H2:H10)
$ws.Range("F2:F10").ClearContents()
$ws.Range("H2:H10").ClearContents()

# Remove the summary statistic formulas (Average, Median, Max, Std. Deviation)
$ws.Range("B15").ClearContents()
$ws.Range("B16").ClearContents()
$ws.Range("B17").ClearContents()
$ws.Range("B19").ClearContents()

# Match the saved sheet view: active cell B15 with B15:B21 selected
$ws.Range("B15:B21").Select()
